$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (rows 2-10), for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
    2  = @{ E=3; G=36.35426266666666; H=109.062788;          I=0.4094848412143908; J=0.4094848412143908; K=3; M=1.564139666666667;  N=4.692419;            O=0.166125338305886;  P=0.166125338305886;  Q=56.86314428935244;  R=511.768298604172;  S=0.06802580777787269; T=0.06802580777787269 }
    3  = @{ E=3; G=36.35426266666666; H=109.062788;          I=0.4094848412143908; J=0.4094848412143908; K=3; M=3.316850333333333;  N=9.950551000000001;   O=0.3522785691569683; P=0.3522785691569683; Q=120.5816482440209;  R=1085.234834196188; S=0.144252733954474;   T=0.144252733954474 }
    4  = @{ E=3; G=36.35426266666666; H=109.062788;          I=0.4094848412143908; J=0.4094848412143908; K=3; M=4.534428999999999;  N=13.603287;            O=0.4815960925371456; P=0.4815960925371456; Q=164.8458229093506;  R=1483.612406184156; S=0.1972062994820442;  T=0.1972062994820442 }
    5  = @{ E=3; G=45.11545066666667; H=135.346352;          I=0.5081685556916724; J=0.5081685556916724; K=3; M=1.564139666666667;  N=4.692419;            O=0.166125338305886;  P=0.166125338305886;  Q=70.56686596727646;  R=635.1017937054881; S=0.08441967323069255; T=0.08441967323069255 }
    6  = @{ E=3; G=45.11545066666667; H=135.346352;          I=0.5081685556916724; J=0.5081685556916724; K=3; M=3.316850333333333;  N=9.950551000000001;   O=0.3522785691569683; P=0.3522785691569683; Q=149.6411975822169;  R=1346.770778239952; S=0.1790168916896255;  T=0.1790168916896255 }
    7  = @{ E=3; G=45.11545066666667; H=135.346352;          I=0.5081685556916724; J=0.5081685556916724; K=3; M=4.534428999999999;  N=13.603287;            O=0.4815960925371456; P=0.4815960925371456; Q=204.5728078510027;  R=1841.155270659024; S=0.2447319907713543;  T=0.2447319907713543 }
    8  = @{ E=3; G=7.310771333333334; H=21.932314;            I=0.0823466030939367; J=0.0823466030939367; K=3; M=1.564139666666667;  N=4.692419;            O=0.166125338305886;  P=0.166125338305886;  Q=11.43506743639622;  R=102.915606927566;  S=0.01367985729732075; T=0.01367985729732075 }
    9  = @{ E=3; G=7.310771333333334; H=21.932314;            I=0.0823466030939367; J=0.0823466030939367; K=3; M=3.316850333333333;  N=9.950551000000001;   O=0.3522785691569683; P=0.3522785691569683; Q=24.24873433389045;  R=218.238609005014;  S=0.0290089435128688;  T=0.0290089435128688 }
    10 = @{ E=3; G=7.310771333333334; H=21.932314;            I=0.0823466030939367; J=0.0823466030939367; K=3; M=4.534428999999999;  N=13.603287;            O=0.4815960925371456; P=0.4815960925371456; Q=33.15017354623533;  R=298.351561916118;  S=0.03965780228374714; T=0.03965780228374714 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
